$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  B=-1.08; C=-36.61; D=665.4580300886175}
    @{Row=3;  B=2.58;  C=2.57;   D=852.2645515146485}
    @{Row=4;  B=1.59;  C=-44.76; D=148.2422796806763}
    @{Row=5;  B=4.5;   C=-23.07; D=159.3633706374936}
    @{Row=6;  B=2.01;  C=87.84999999999999; D=712.7101401789072}
    @{Row=7;  B=1.21;  C=-27.84; D=406.7558751760207}
    @{Row=8;  B=-2.73; C=-37.17; D=233.984900923297}
    @{Row=9;  B=-4.5;  C=-74.67; D=310.5968891251165}
    @{Row=10; B=1.7;   C=6.65;   D=645}
    @{Row=11; B=0.67;  C=-7.25;  D=434.3055113215779}
    @{Row=12; B=-1.69; C=44.31;  D=817.257859779848}
    @{Row=13; B=5.72;  C=-19.32; D=877.9071361015275}
    @{Row=14; B=0.82;  C=29.07;  D=718.3415849059797}
    @{Row=15; B=0.43;  C=-24.5;  D=441.0003016115479}
    @{Row=16; B=3.25;  C=32.91;  D=548.8453146526005}
    @{Row=17; B=4.83;  C=-75.68000000000001; D=69.05455730801823}
    @{Row=18; B=1.31;  C=-19.59; D=8041.088303005878}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
}
